# Update the date heading (text is unique in the document, so a normal
# document-wide Find/Replace is safe here).
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-07-25 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-26 Friday", 2) | Out-Null

# Update the 5x5 block of division problems held in the single table.
# Data rows are 1, 5, 9, 13, 17 (1-indexed); each has 5 populated cells.
#
# NOTE: Find.Execute on a Cell's Range does NOT respect that range's end
# boundary in this runtime -- it keeps searching forward through the rest
# of the document/story, so it can "find" and clobber text that lives in a
# later cell (this matters because some old/new values coincide, e.g.
# 60÷8= -> 96÷5= while a separate original cell already reads 96÷5= ->
# 16÷6=). To keep each edit strictly confined to its own cell, assign the
# new text directly to the cell's Range instead of using Find/Replace.
# Assigning Range.Text only swaps the run's text content and leaves the
# run's formatting (rFonts/sz/etc.) untouched.

$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $table.Cell($row, $col).Range.Text = $newText
}

Set-CellText $t 1 1 "83÷2="
Set-CellText $t 1 2 "92÷3="
Set-CellText $t 1 3 "20÷6="
Set-CellText $t 1 4 "34÷3="
Set-CellText $t 1 5 "72÷6="

Set-CellText $t 5 1 "61÷2="
Set-CellText $t 5 2 "50÷3="
Set-CellText $t 5 3 "37÷6="
Set-CellText $t 5 4 "26÷7="
Set-CellText $t 5 5 "96÷5="

Set-CellText $t 9 1 "16÷6="
Set-CellText $t 9 2 "28÷3="
Set-CellText $t 9 3 "28÷2="
Set-CellText $t 9 4 "52÷8="
Set-CellText $t 9 5 "34÷6="

Set-CellText $t 13 1 "43÷9="
Set-CellText $t 13 2 "97÷4="
Set-CellText $t 13 3 "50÷3="
Set-CellText $t 13 4 "26÷4="
Set-CellText $t 13 5 "41÷8="

Set-CellText $t 17 1 "36÷3="
Set-CellText $t 17 2 "15÷7="
Set-CellText $t 17 3 "20÷7="
Set-CellText $t 17 4 "92÷6="
Set-CellText $t 17 5 "68÷3="
